{"js": "// Change 1: remove the trailing sentence about assuming a maximum string\n// length, which followed \"...is wrapper over a hash table.\" in the\n// \"Insert\" method paragraph. The sentence (including its leading space)\n// is deleted outright, leaving the period after \"hash table\" intact.\nconst removeResults = context.document.body.search(\n  \" For the rest of the analysis, the string length will be assumed to have a maximum, which results in the insert function being treated as running in constant time.\",\n  { matchCase: true }\n);\nremoveResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < removeResults.items.length; i++) {\n  removeResults.items[i].insertText(\"\", \"Replace\");\n}\nawait context.sync();\n\n// Change 2: rewrite the final \"time complexity of the main method\" summary\n// paragraph to talk about the Page Rank algorithm's complexity instead of\n// summing the three PageRank class methods.\nconst oldSummary =\n  \"The time complexity of the main method is the sum of the three methods in the PageRank class; O(|E|) + O(p * |E|) + O(|V|). Combining and eliminating constant factors results in the entire project having a runtime complexity of O(|E| + |V|).\";\nconst newSummary =\n  \"The time complexity of the main Page Rank algorithm is O(p * |E|) as discussed above. If we consider p to be a constant, this reduces down to a time complexity of O(|E|).\";\n\nconst summaryResults = context.document.body.search(oldSummary, { matchCase: true });\nsummaryResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < summaryResults.items.length; i++) {\n  summaryResults.items[i].insertText(newSummary, \"Replace\");\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Change 1: remove the trailing sentence about assuming a maximum string\n# length, which followed \"...is wrapper over a hash table.\" in the\n# \"Insert\" method paragraph. The sentence (including its leading space)\n# is deleted outright, leaving the period after \"hash table\" intact.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\n    \" For the rest of the analysis, the string length will be assumed to have a maximum, which results in the insert function being treated as running in constant time.\",\n    $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2\n) | Out-Null\n\n# Change 2: rewrite the final \"time complexity of the main method\" summary\n# paragraph to talk about the Page Rank algorithm's complexity instead of\n# summing the three PageRank class methods.\n$oldSummary = \"The time complexity of the main method is the sum of the three methods in the PageRank class; O(|E|) + O(p * |E|) + O(|V|). Combining and eliminating constant factors results in the entire project having a runtime complexity of O(|E| + |V|).\"\n$newSummary = \"The time complexity of the main Page Rank algorithm is O(p * |E|) as discussed above. If we consider p to be a constant, this reduces down to a time complexity of O(|E|).\"\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\n    $oldSummary,\n    $false, $false, $false, $false, $false, $true, 1, $false, $newSummary, 2\n) | Out-Null\n"}
